# fix: Fixed column values
# Add a new worksheet "Variacion - Retorno" after "Tenencia" that mirrors
# the existing sheet's layout but with corrected (zeroed) return values
# in column B, and make the new sheet the active tab.

$wb = $excel.ActiveWorkbook
$tenencia = $wb.Worksheets.Item(1)

$new = $wb.Worksheets.Add($null, $tenencia)
$new.Name = "Variacion - Retorno"

# Header (merged cell)
$new.Range("B1").Value = "TASA FIJA"
$new.Range("B1:B1").Merge()

# Column headers
$new.Range("A2").Value = "DateRequested"
$new.Range("B2").Value = "S31E5"

# Data rows - dates kept as before, values fixed to 0.000000
$new.Range("A3").Value = "'2024-08-01"
$new.Range("A3").Style = "Normal"
$new.Range("B3").Value = "'0.000000"
$new.Range("B3").Style = "Normal"

$new.Range("A4").Value = "'2024-08-02"
$new.Range("A4").Style = "Normal"
$new.Range("B4").Value = "'0.000000"
$new.Range("B4").Style = "Normal"

$new.Range("A5").Value = "'2024-08-03"
$new.Range("A5").Style = "Normal"
$new.Range("B5").Value = "'0.000000"
$new.Range("B5").Style = "Normal"

# Make the new sheet the active tab (also clears tabSelected on "Tenencia")
$new.Activate()
